$d = $word.ActiveDocument

$d.Content.Find.Execute("668×6=", $true, $false, $false, $false, $false, $true, 1, $false, "430×6=", 2) | Out-Null
$d.Content.Find.Execute("189×5=", $true, $false, $false, $false, $false, $true, 1, $false, "891×9=", 2) | Out-Null
$d.Content.Find.Execute("259×9=", $true, $false, $false, $false, $false, $true, 1, $false, "813×3=", 2) | Out-Null
$d.Content.Find.Execute("208×7=", $true, $false, $false, $false, $false, $true, 1, $false, "961×7=", 2) | Out-Null
$d.Content.Find.Execute("529×9=", $true, $false, $false, $false, $false, $true, 1, $false, "878×8=", 2) | Out-Null
$d.Content.Find.Execute("914×4=", $true, $false, $false, $false, $false, $true, 1, $false, "302×2=", 2) | Out-Null
$d.Content.Find.Execute("858×4=", $true, $false, $false, $false, $false, $true, 1, $false, "109×9=", 2) | Out-Null
$d.Content.Find.Execute("485×6=", $true, $false, $false, $false, $false, $true, 1, $false, "783×3=", 2) | Out-Null
$d.Content.Find.Execute("161×4=", $true, $false, $false, $false, $false, $true, 1, $false, "689×9=", 2) | Out-Null
$d.Content.Find.Execute("541×5=", $true, $false, $false, $false, $false, $true, 1, $false, "582×4=", 2) | Out-Null
$d.Content.Find.Execute("887×9=", $true, $false, $false, $false, $false, $true, 1, $false, "905×8=", 2) | Out-Null
$d.Content.Find.Execute("729×6=", $true, $false, $false, $false, $false, $true, 1, $false, "174×2=", 2) | Out-Null
$d.Content.Find.Execute("695×6=", $true, $false, $false, $false, $false, $true, 1, $false, "293×7=", 2) | Out-Null
$d.Content.Find.Execute("943×7=", $true, $false, $false, $false, $false, $true, 1, $false, "970×8=", 2) | Out-Null
$d.Content.Find.Execute("874×6=", $true, $false, $false, $false, $false, $true, 1, $false, "408×3=", 2) | Out-Null
$d.Content.Find.Execute("179×8=", $true, $false, $false, $false, $false, $true, 1, $false, "504×9=", 2) | Out-Null
$d.Content.Find.Execute("274×3=", $true, $false, $false, $false, $false, $true, 1, $false, "657×4=", 2) | Out-Null
$d.Content.Find.Execute("883×6=", $true, $false, $false, $false, $false, $true, 1, $false, "992×9=", 2) | Out-Null
$d.Content.Find.Execute("829×3=", $true, $false, $false, $false, $false, $true, 1, $false, "349×3=", 2) | Out-Null
$d.Content.Find.Execute("965×7=", $true, $false, $false, $false, $false, $true, 1, $false, "798×7=", 2) | Out-Null
$d.Content.Find.Execute("992×5=", $true, $false, $false, $false, $false, $true, 1, $false, "430×3=", 2) | Out-Null
$d.Content.Find.Execute("182×2=", $true, $false, $false, $false, $false, $true, 1, $false, "742×7=", 2) | Out-Null
$d.Content.Find.Execute("712×5=", $true, $false, $false, $false, $false, $true, 1, $false, "422×5=", 2) | Out-Null
$d.Content.Find.Execute("954×4=", $true, $false, $false, $false, $false, $true, 1, $false, "153×3=", 2) | Out-Null
$d.Content.Find.Execute("566×6=", $true, $false, $false, $false, $false, $true, 1, $false, "157×3=", 2) | Out-Null
